$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 189 ---
$ws.Range("A189").Value = 617
$ws.Range("C189").Value = 29
$ws.Range("D189").Value = 'Female'
$ws.Range("E189").Value = 'Y'
$ws.Range("F189").Value = 'N'
$ws.Range("G189").Value = 'N'
$ws.Range("H189").Value = 'UK'
$ws.Range("I189").Value = 'UK'
$ws.Range("J189").Value = 'Y'
$ws.Range("K189").Value = 'N'
$ws.Range("L189").Value = 'UK'
$ws.Range("M189").Value = 'N'
$ws.Range("N189").Value = 'UK'
$ws.Range("O189").Value = 'Y'
$ws.Range("P189").Value = 'N'
$ws.Range("Q189").Value = 'N'
$ws.Range("R189").Value = 'N'
$ws.Range("T189").Value = 'N'
$ws.Range("U189").Value = 'N'
$ws.Range("V189").Value = 'N'
$ws.Range("W189").Value = 'N'
$ws.Range("X189").Value = 'Y'
$ws.Range("Y189").Value = 'Y'
$ws.Range("Z189").Value = '"ten razor blades after coating them with chewing gum"'
$ws.Range("S189").Value = '"history of psychiatric disorder that was treated with thorazine and lithium, but was not taking any of her medications", "denied any intentions to hurt herself"'
$ws.Range("AA189").Value = 'N'
$ws.Range("AB189").Value = 'N'
$ws.Range("AC189").Value = 'N'
$ws.Range("AD189").Value = 'N'
$ws.Range("AE189").Value = 'N'
$ws.Range("AF189").Value = 'N'
$ws.Range("AG189").Value = 'N'
$ws.Range("AH189").Value = '"Eventually, all razor blades were passed"'
$ws.Range("B189").Value = '617-001'

# --- Row 190 ---
$ws.Range("A190").Value = 620
$ws.Range("B190").Value = '620-001'
$ws.Range("C190").Value = 32
$ws.Range("D190").Value = 'Female'
$ws.Range("E190").Value = 'Y'
$ws.Range("F190").Value = 'N'
$ws.Range("G190").Value = 'N'
$ws.Range("H190").Value = 'UK'
$ws.Range("I190").Value = 'UK'
$ws.Range("J190").Value = 'Y'
$ws.Range("K190").Value = 'N'
$ws.Range("L190").Value = 'Y'
$ws.Range("M190").Value = 'UK'
$ws.Range("N190").Value = 'UK'
$ws.Range("O190").Value = 'UK'
$ws.Range("P190").Value = 'UK'
$ws.Range("Q190").Value = 'Y'
$ws.Range("R190").Value = 'N'
$ws.Range("S190").Value = '"psychiatric history of compulsive foreign body ingestions resulting in multiple upper endoscopies and retrievals"'
$ws.Range("T190").Value = 'N'
$ws.Range("U190").Value = 'N'
$ws.Range("V190").Value = 'Y'
$ws.Range("W190").Value = 'N'
$ws.Range("X190").Value = 'Y'
$ws.Range("Y190").Value = 'N'
$ws.Range("Z190").Value = '"butter knife"'
$ws.Range("AA190").Value = 'Y'
$ws.Range("AB190").Value = 'Y'
$ws.Range("AC190").Value = 'N'
$ws.Range("AD190").Value = 'Y'
$ws.Range("AE190").Value = 'Y'
$ws.Range("AF190").Value = 'N'
$ws.Range("AG190").Value = 'N'

# --- Row 191 ---
$ws.Range("A191").Value = 621
$ws.Range("C191").Value = 24
$ws.Range("D191").Value = 'Male'
$ws.Range("E191").Value = 'Y'
$ws.Range("F191").Value = 'N'
$ws.Range("G191").Value = 'N'
$ws.Range("H191").Value = 'UK'
$ws.Range("I191").Value = 'UK'
$ws.Range("J191").Value = 'Y'
$ws.Range("K191").Value = 'N'
$ws.Range("L191").Value = 'UK'
$ws.Range("M191").Value = 'UK'
$ws.Range("N191").Value = 'UK'
$ws.Range("O191").Value = 'UK'
$ws.Range("P191").Value = 'UK'
$ws.Range("Q191").Value = 'Y'
$ws.Range("R191").Value = 'N'
$ws.Range("S191").Value = '"Maniac Depressive Psychosis who had a history of ingesting Nails and screws of sizes varying from 2 cm to 15 cm for more than 1 year without causing any perforation and other acute complication"'
$ws.Range("B191").Value = '621-001'
$ws.Range("T191").Value = 'N'
$ws.Range("U191").Value = 'N'
$ws.Range("V191").Value = 'Y'
$ws.Range("W191").Value = 'Y'
$ws.Range("X191").Value = 'Y'
$ws.Range("Y191").Value = 'Y'
$ws.Range("Z191").Value = '"27 metal nails and screws of sizes 6 cm to 15 cm and bent in various shapes were removed from inside the stomach."'
$ws.Range("AA191").Value = 'N'
$ws.Range("AB191").Value = 'Y'
$ws.Range("AC191").Value = 'N'
$ws.Range("AD191").Value = 'N'
$ws.Range("AE191").Value = 'N'
$ws.Range("AF191").Value = 'N'
$ws.Range("AG191").Value = 'N'

# --- Row 192 ---
$ws.Range("A192").Value = 622
$ws.Range("B192").Value = '622-001'
$ws.Range("C192").Value = 100
$ws.Range("D192").Value = 'Female'
$ws.Range("E192").Value = 'Y'
$ws.Range("F192").Value = 'N'
$ws.Range("G192").Value = 'N'
$ws.Range("H192").Value = 'UK'
$ws.Range("I192").Value = 'UK'
$ws.Range("J192").Value = 'N'
$ws.Range("K192").Value = 'N'
$ws.Range("L192").Value = 'N'
$ws.Range("M192").Value = 'Y'
$ws.Range("N192").Value = 'N'
$ws.Range("O192").Value = 'Y'
$ws.Range("P192").Value = 'Y'
$ws.Range("Q192").Value = 'N'
$ws.Range("R192").Value = 'N'
$ws.Range("T192").Value = 'N'
$ws.Range("U192").Value = 'N'
$ws.Range("V192").Value = 'N'
$ws.Range("W192").Value = 'N'
$ws.Range("X192").Value = 'N'
$ws.Range("Y192").Value = 'Y'
$ws.Range("Z192").Value = '"26 coins, one ferrous ring and one cylindrical plastic object were retrieved"'
$ws.Range("S192").Value = '"suicide attempt due to intolerable pain induced by a fracture she suffered 3 mo previously. She was bedridden and had a depressed mood, which caused her to attempt suicide"'
$ws.Range("AA192").Value = 'Y'
$ws.Range("AB192").Value = 'N'
$ws.Range("AC192").Value = 'N'
$ws.Range("AD192").Value = 'N'
$ws.Range("AE192").Value = 'N'
$ws.Range("AF192").Value = 'N'
$ws.Range("AG192").Value = 'N'

# --- Row 193 ---
$ws.Range("A193").Value = 623
$ws.Range("B193").Value = '623-001'
$ws.Range("C193").Value = 26
$ws.Range("D193").Value = 'Male'
$ws.Range("E193").Value = 'Y'
$ws.Range("F193").Value = 'N'
$ws.Range("G193").Value = 'UK'
$ws.Range("H193").Value = 'UK'
$ws.Range("I193").Value = 'UK'
$ws.Range("J193").Value = 'Y'
$ws.Range("K193").Value = 'N'
$ws.Range("L193").Value = 'UK'
$ws.Range("M193").Value = 'UK'
$ws.Range("N193").Value = 'UK'
$ws.Range("O193").Value = 'Y'
$ws.Range("P193").Value = 'UK'
$ws.Range("Q193").Value = 'N'
$ws.Range("R193").Value = 'N'
$ws.Range("S193").Value = '"known psychiatric illness", "ingesting sewing needles by wrapping it on a plant leaf out of a schizophrenic disorder"'
$ws.Range("T193").Value = 'N'
$ws.Range("U193").Value = 'N'
$ws.Range("V193").Value = 'Y'
$ws.Range("W193").Value = 'Y'
$ws.Range("X193").Value = 'Y'
$ws.Range("Y193").Value = 'Y'
$ws.Range("Z193").Value = '"ingesting sewing needles by wrapping it on a plant leaf out of a schizophrenic disorder", "8cm long"'
$ws.Range("AA193").Value = 'N'
$ws.Range("AB193").Value = 'Y'
$ws.Range("AC193").Value = 'N'
$ws.Range("AD193").Value = 'Y'
$ws.Range("AE193").Value = 'Y'
$ws.Range("AF193").Value = 'N'
$ws.Range("AG193").Value = 'Y'
$ws.Range("AH193").Value = '"perforation peritonitis", "needle in liver", "needle in ballder"'

# --- Row 194 ---
$ws.Range("A194").Value = 625
$ws.Range("B194").Value = '625-001'
$ws.Range("C194").Value = 30
$ws.Range("D194").Value = 'Female'
$ws.Range("E194").Value = 'Y'
$ws.Range("F194").Value = 'N'
$ws.Range("G194").Value = 'N'
$ws.Range("H194").Value = 'UK'
$ws.Range("I194").Value = 'UK'
$ws.Range("J194").Value = 'Y'
$ws.Range("K194").Value = 'N'
$ws.Range("L194").Value = 'Y'
$ws.Range("M194").Value = 'UK'
$ws.Range("N194").Value = 'UK'
$ws.Range("O194").Value = 'UK'
$ws.Range("P194").Value = 'UK'
$ws.Range("Q194").Value = 'Y'
$ws.Range("R194").Value = 'N'
$ws.Range("S194").Value = '"medical history of borderline personality disorder, atypical schizo-affective disorder, a long history of self-mutilating behavior and numerous suicide attempts"'
$ws.Range("T194").Value = 'N'
$ws.Range("U194").Value = 'N'
$ws.Range("V194").Value = 'Y'
$ws.Range("W194").Value = 'Y'
$ws.Range("X194").Value = 'Y'
$ws.Range("Y194").Value = 'Y'
$ws.Range("Z194").Value = '"10 razor blades wrapped in paper and chewing gum"'
$ws.Range("AA194").Value = 'Y'
$ws.Range("AB194").Value = 'N'
$ws.Range("AC194").Value = 'N'
$ws.Range("AD194").Value = 'N'
$ws.Range("AE194").Value = 'N'
$ws.Range("AF194").Value = 'N'
$ws.Range("AG194").Value = 'Y'
$ws.Range("AH194").Value = '"she passed the razor blades transanally without further event"'

# Update selection to match final cursor position
$ws.Range("R194").Select()

